# Reorders same-matchday rows into ascending id (column B) order.
# The scraped feed wrote a few matchdays with rows out of id sequence; this
# swaps the full match record (columns B and E:AB) between rows while leaving
# column A (the sequential running index) untouched, matching the source fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32, 33
$row32 = $ws.Range("B32:AB32").Value2
$row33 = $ws.Range("B33:AB33").Value2
$ws.Range("B32:AB32").Value2 = $row33
$ws.Range("B33:AB33").Value2 = $row32

# Rows 60, 61
$row60 = $ws.Range("B60:AB60").Value2
$row61 = $ws.Range("B61:AB61").Value2
$ws.Range("B60:AB60").Value2 = $row61
$ws.Range("B61:AB61").Value2 = $row60

# Rows 140, 141
$row140 = $ws.Range("B140:AB140").Value2
$row141 = $ws.Range("B141:AB141").Value2
$ws.Range("B140:AB140").Value2 = $row141
$ws.Range("B141:AB141").Value2 = $row140

# Rows 142, 143
$row142 = $ws.Range("B142:AB142").Value2
$row143 = $ws.Range("B143:AB143").Value2
$ws.Range("B142:AB142").Value2 = $row143
$ws.Range("B143:AB143").Value2 = $row142

# Rows 151, 152
$row151 = $ws.Range("B151:AB151").Value2
$row152 = $ws.Range("B152:AB152").Value2
$ws.Range("B151:AB151").Value2 = $row152
$ws.Range("B152:AB152").Value2 = $row151

# Rows 167, 168
$row167 = $ws.Range("B167:AB167").Value2
$row168 = $ws.Range("B168:AB168").Value2
$ws.Range("B167:AB167").Value2 = $row168
$ws.Range("B168:AB168").Value2 = $row167

# Rows 186, 187
$row186 = $ws.Range("B186:AB186").Value2
$row187 = $ws.Range("B187:AB187").Value2
$ws.Range("B186:AB186").Value2 = $row187
$ws.Range("B187:AB187").Value2 = $row186

# Rows 201, 202
$row201 = $ws.Range("B201:AB201").Value2
$row202 = $ws.Range("B202:AB202").Value2
$ws.Range("B201:AB201").Value2 = $row202
$ws.Range("B202:AB202").Value2 = $row201

# Rows 221, 222, 223
$row221 = $ws.Range("B221:AB221").Value2
$row222 = $ws.Range("B222:AB222").Value2
$row223 = $ws.Range("B223:AB223").Value2
$ws.Range("B221:AB221").Value2 = $row223
$ws.Range("B222:AB222").Value2 = $row221
$ws.Range("B223:AB223").Value2 = $row222

# Rows 230, 231
$row230 = $ws.Range("B230:AB230").Value2
$row231 = $ws.Range("B231:AB231").Value2
$ws.Range("B230:AB230").Value2 = $row231
$ws.Range("B231:AB231").Value2 = $row230

# Rows 241, 242
$row241 = $ws.Range("B241:AB241").Value2
$row242 = $ws.Range("B242:AB242").Value2
$ws.Range("B241:AB241").Value2 = $row242
$ws.Range("B242:AB242").Value2 = $row241

# Rows 260, 261
$row260 = $ws.Range("B260:AB260").Value2
$row261 = $ws.Range("B261:AB261").Value2
$ws.Range("B260:AB260").Value2 = $row261
$ws.Range("B261:AB261").Value2 = $row260

# Rows 278, 279
$row278 = $ws.Range("B278:AB278").Value2
$row279 = $ws.Range("B279:AB279").Value2
$ws.Range("B278:AB278").Value2 = $row279
$ws.Range("B279:AB279").Value2 = $row278

# Rows 293, 294
$row293 = $ws.Range("B293:AB293").Value2
$row294 = $ws.Range("B294:AB294").Value2
$ws.Range("B293:AB293").Value2 = $row294
$ws.Range("B294:AB294").Value2 = $row293

# Rows 296, 297
$row296 = $ws.Range("B296:AB296").Value2
$row297 = $ws.Range("B297:AB297").Value2
$ws.Range("B296:AB296").Value2 = $row297
$ws.Range("B297:AB297").Value2 = $row296

# Rows 303, 304
$row303 = $ws.Range("B303:AB303").Value2
$row304 = $ws.Range("B304:AB304").Value2
$ws.Range("B303:AB303").Value2 = $row304
$ws.Range("B304:AB304").Value2 = $row303

# Rows 306, 307
$row306 = $ws.Range("B306:AB306").Value2
$row307 = $ws.Range("B307:AB307").Value2
$ws.Range("B306:AB306").Value2 = $row307
$ws.Range("B307:AB307").Value2 = $row306
